$d = $word.ActiveDocument

# 1. Remove the comment entirely. This also strips the commentRangeStart /
#    commentRangeEnd / commentReference markers from document.xml and empties
#    comments.xml / commentsExtended.xml / commentsExtensible.xml /
#    commentsIds.xml / people.xml.
if ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# 2. Rework the sentence:
#      "These data were published as part of the Healthy Rivers and
#       Landscapes program."
#    into:
#      "These data were published to support the Healthy Rivers and
#       Landscapes Program."
#    split across separate runs, mirroring the authored edit.

$oldPhrase = "These data were published as part of the Healthy Rivers and Landscapes program."
$full = $d.Content.Text
$start = $full.IndexOf($oldPhrase)

$prefix = "These data were published "

# -- Move the run boundary forward so "These data were published " sticks to
#    the preceding sentence run instead of starting its own run. We force a
#    run split at this boundary by toggling formatting on/off (back to the
#    original look) so the range stays its own <w:r> going forward.
$boundary = $d.Range($start, $start + $prefix.Length)
$boundary.Bold = 1
$boundary.Bold = 0

$pos = $start + $prefix.Length

# -- "as part of" -> "to support"
$newText = "to support"
$r = $d.Range($pos, $pos + "as part of".Length)
$r.Text = $newText
$r.Bold = 1
$r.Bold = 0
$pos = $pos + $newText.Length

# -- " the Healthy Rivers and Landscapes " (unchanged text, own run)
$newText = " the Healthy Rivers and Landscapes "
$r = $d.Range($pos, $pos + $newText.Length)
$r.Text = $newText
$r.Bold = 1
$r.Bold = 0
$pos = $pos + $newText.Length

# -- "p" -> "P"
$newText = "P"
$r = $d.Range($pos, $pos + "p".Length)
$r.Text = $newText
$r.Bold = 1
$r.Bold = 0
$pos = $pos + $newText.Length

# -- "rogram." (unchanged text, own run)
$newText = "rogram."
$r = $d.Range($pos, $pos + $newText.Length)
$r.Text = $newText
$pos = $pos + $newText.Length

Write-Host "Result: " $d.Content.Text.Substring($start - 5, 150)
